# Updated cryptos list on Tue Nov 21 08:59:02 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.441.37"
$ws.Range("E2").Value = "  +0.83%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.010.68"
$ws.Range("E3").Value = "  -0.28%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "260.26"
$ws.Range("E5").Value = "  +5.34%  "
$ws.Range("E6").Value = "  -0.75%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "56.35"
$ws.Range("E8").Value = "  -5.70%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.385"
$ws.Range("E9").Value = "  -1.26%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0772"
$ws.Range("E10").Value = "  -4.39%  "
$ws.Range("E11").Value = "  -2.72%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.28"
$ws.Range("E12").Value = "  -5.92%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.304.38"
$ws.Range("E13").Value = "  -0.32%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.05"
$ws.Range("E14").Value = "  -5.61%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.801"
$ws.Range("E15").Value = "  -5.71%  "
$ws.Range("E16").Value = "  -4.67%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.000.31"
$ws.Range("E17").Value = "  -1.37%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "37.337.39"
$ws.Range("E18").Value = "  +0.81%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "69.81"
$ws.Range("E19").Value = "  -0.92%  "
$ws.Range("E20").Value = "  -3.46%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.15"
$ws.Range("E21").Value = "  -1.35%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "231.77"
$ws.Range("E22").Value = "  +0.62%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.62"
$ws.Range("E23").Value = "  +4.40%  "
$ws.Range("E24").Value = "  +0.05%  "
$ws.Range("E25").Value = "  -0.56%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "164.59"
$ws.Range("E26").Value = "  +0.63%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.90"
$ws.Range("E27").Value = "  -5.53%  "
$ws.Range("E28").Value = "  -0.48%  "
$ws.Range("E29").Value = "  -4.26%  "
$ws.Range("E30").Value = "  -2.14%  "
$ws.Range("E31").Value = "  -1.52%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.62"
$ws.Range("E32").Value = "  -3.97%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0646"
$ws.Range("E33").Value = "  -2.14%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.53"
$ws.Range("E34").Value = "  +0.74%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.37"
$ws.Range("E35").Value = "  -3.30%  "
$ws.Range("E36").Value = "  +0.41%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.36"
$ws.Range("E38").Value = "  -4.66%  "
$ws.Range("E39").Value = "  -0.46%  "
$ws.Range("E40").Value = "  +3.91%  "
$ws.Range("E41").Value = "  +0.25%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0213"
$ws.Range("E42").Value = "  -0.85%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0925"
$ws.Range("E43").Value = "  -6.18%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.414.17"
$ws.Range("E44").Value = "  +2.00%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "15.81"
$ws.Range("E45").Value = "  -5.12%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "89.54"
$ws.Range("E46").Value = "  -2.56%  "
$ws.Range("E47").Value = "  -3.25%  "
$ws.Range("E48").Value = "  +2.32%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.04"
$ws.Range("E49").Value = "  -6.01%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.195.83"
$ws.Range("E50").Value = "  -0.36%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.95"
$ws.Range("E51").Value = "  -9.22%  "
